$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need to be forced to
# Text format first, otherwise Excel auto-converts them to numbers.
$textCells = @("D5", "D7", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D30", "D31", "D33", "D36", "D38", "D39", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.433.89'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '2.243.88'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '244.99'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("E6").Value = '  +1.05%  '
$ws.Range("D7").Value = '75.68'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").Value = '44.33'
$ws.Range("E10").Value = '  +8.54%  '
$ws.Range("D11").Value = '0.0949'
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").Value = '7.23'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").Value = '14.60'
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("D15").Value = '0.860'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").Value = '2.251.21'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '42.301.64'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("E18").Value = '  +4.03%  '
$ws.Range("D19").Value = '6.20'
$ws.Range("E19").Value = '  +1.03%  '
$ws.Range("D20").Value = '72.04'
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").Value = '11.38'
$ws.Range("E21").Value = '  +58.17%  '
$ws.Range("D22").Value = '2.24'
$ws.Range("E22").Value = '  -3.42%  '
$ws.Range("D23").Value = '231.78'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '11.83'
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -1.43%  '
$ws.Range("D27").Value = '2.30'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  +4.97%  '
$ws.Range("D29").Value = '167.13'
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("D30").Value = '20.72'
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("D31").Value = '5.82'
$ws.Range("E31").Value = '  +19.42%  '
$ws.Range("E32").Value = '  -2.13%  '
$ws.Range("D33").Value = '30.93'
$ws.Range("E33").Value = '  -5.77%  '
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").Value = '4.69'
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("E37").Value = '  +4.64%  '
$ws.Range("D38").Value = '13.86'
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("D39").Value = '2.18'
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").Value = '5.77'
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("D41").Value = '63.76'
$ws.Range("E41").Value = '  +4.30%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '106.95'
$ws.Range("E43").Value = '  -5.29%  '
$ws.Range("D44").Value = '8.85'
$ws.Range("E44").Value = '  +1.49%  '
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("D46").Value = '0.996'
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").Value = '2.44'
$ws.Range("E47").Value = '  +7.76%  '
$ws.Range("D48").Value = '1.15'
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = '2.72'
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = '4.13'
$ws.Range("E51").Value = '  -1.08%  '
